$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 357
$ws1.Cells.Item(3, 6).Value = 282
$ws1.Cells.Item(4, 6).Value = 1809
$ws1.Cells.Item(6, 6).Value = 30
$ws1.Cells.Item(7, 6).Value = 62
$ws1.Cells.Item(9, 6).Value = 153
$ws1.Cells.Item(10, 6).Value = 3567
$ws1.Cells.Item(12, 6).Value = 94
$ws1.Cells.Item(15, 6).Value = 60
$ws1.Cells.Item(16, 6).Value = 621
$ws1.Cells.Item(17, 6).Value = 114
$ws1.Cells.Item(18, 6).Value = 782
$ws1.Cells.Item(19, 6).Value = 7
$ws1.Cells.Item(20, 6).Value = 216
$ws1.Cells.Item(21, 6).Value = 135
$ws1.Cells.Item(23, 6).Value = 67
$ws1.Cells.Item(25, 6).Value = 2795
$ws1.Cells.Item(26, 6).Value = 5272
$ws1.Cells.Item(30, 6).Value = 3097
$ws1.Cells.Item(32, 6).Value = 2283
$ws1.Cells.Item(35, 6).Value = 89
$ws1.Cells.Item(37, 6).Value = 191
$ws1.Cells.Item(42, 6).Value = 31
$ws1.Cells.Item(43, 6).Value = 17
$ws1.Cells.Item(44, 6).Value = 455
$ws1.Cells.Item(46, 6).Value = 500

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 357
$ws4.Cells.Item(3, 6).Value = 282
$ws4.Cells.Item(4, 6).Value = 1809
$ws4.Cells.Item(6, 6).Value = 30
$ws4.Cells.Item(7, 6).Value = 62
$ws4.Cells.Item(9, 6).Value = 153
$ws4.Cells.Item(10, 6).Value = 3567
$ws4.Cells.Item(12, 6).Value = 94
$ws4.Cells.Item(16, 6).Value = 60
$ws4.Cells.Item(17, 6).Value = 621
$ws4.Cells.Item(18, 6).Value = 114
$ws4.Cells.Item(19, 6).Value = 782
$ws4.Cells.Item(20, 6).Value = 7
$ws4.Cells.Item(21, 6).Value = 216
$ws4.Cells.Item(22, 6).Value = 135
$ws4.Cells.Item(24, 6).Value = 67
$ws4.Cells.Item(26, 6).Value = 2796
$ws4.Cells.Item(27, 6).Value = 5272
$ws4.Cells.Item(31, 6).Value = 3097
$ws4.Cells.Item(33, 6).Value = 2283
$ws4.Cells.Item(36, 6).Value = 89
$ws4.Cells.Item(38, 6).Value = 191
$ws4.Cells.Item(43, 6).Value = 31
$ws4.Cells.Item(44, 6).Value = 17
$ws4.Cells.Item(45, 6).Value = 455
$ws4.Cells.Item(47, 6).Value = 500
